# B6-PowerPoint.pptx edit: Sun, Apr 26, 2020 12:05:06 PM
#
# 1) Re-style the three tables (slides 14, 15, 16) from table style
#    {C695E72F-C4D5-4FCB-929C-FABF501061D2} to {D2F4840E-1226-4A6C-8AFE-39D02CC90317}.
# 2) Recolour the deck's theme so the slide-master theme carries the
#    "Office" colour scheme (dk2/lt2/accent1-6/hlink/folHlink) instead of
#    the previous "Red Violet"/Integral colours - the colour values that
#    the notes-master theme (theme2.xml) already used.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$tableSlides = 14, 15, 16
$newStyleId = "{D2F4840E-1226-4A6C-8AFE-39D02CC90317}"

foreach ($slideIndex in $tableSlides) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# RGB() values below are the standard COM 0x00BBGGRR packing of the
# "Office" theme colours (dk1 000000, lt1 FFFFFF, dk2 44546A, lt2 E7E6E6,
# accent1 5B9BD5, accent2 ED7D31, accent3 A5A5A5, accent4 FFC000,
# accent5 4472C4, accent6 70AD47, hlink 0563C1, folHlink 954F72).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
